$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "D0.jpg"
$ws.Range("A3").Value = "D1.jpg"
$ws.Range("A4").Value = "D2.jpg"
$ws.Range("A5").Value = "D3.jpg"
$ws.Range("A6").Value = "D4.jpg"
$ws.Range("A7").Value = "L0.jpg"
$ws.Range("A8").Value = "L1.jpg"
$ws.Range("A9").Value = "L2.jpg"
$ws.Range("A10").Value = "L3.jpg"
$ws.Range("A11").Value = "L4.jpg"

$ws.Range("K14").Select()
